$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.003642638115616127
$ws.Range("B3").Value = 0.7577246201537545
$ws.Range("B4").Value = 0.1005667238203967
$ws.Range("B5").Value = "[1.0, 0.05202639084398948, -0.04787658854942903, -0.137564655772672, -0.1135759129650511, 0.0814652442430326, 0.08392292471032542, 0.0650141108321826, 0.01882596158466086, -0.038275124714904935, -0.055957743162891295, -0.03351964838726925, -0.027860132598025427, 0.03596356767971705, 0.06635738590000892, 0.02904448965185704, -0.019894208338746715, -0.05739439196321059, -0.07880625685022305, -0.03848955196939101]"

$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.174468039891417
$ws.Range("B5").Value = 0.5526466161711083
$ws.Range("B6").Value = 0.05884290129237798
$ws.Range("B7").Value = 0.1508453391109807
$ws.Range("B8").Value = 0.6086437520559146
$ws.Range("B9").Value = 0.05513269080032153
$ws.Range("B10").Value = "[0.9999999999999998, 0.07910463187441695, 0.0012200437325331835, -0.015339275560061299, 0.012759060733731663, 0.04396598501313142, 0.022868912824314865, 0.009197891318008416, 0.014855828924272337, -0.03592794433053699, 0.017593264928896724, -0.003175085595637512, -0.015624586572985771, 0.02107159620831528, 0.01032916853707136, -0.004872857247384475, -0.01724781518624857, -0.0037221460144831347, -0.00014745470997201595, -0.035560880430825365]"

$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.001267570511828704
$ws.Range("B3").Value = 0.1027706486870374
$ws.Range("B4").Value = 0.06267436121305359
$ws.Range("B5").Value = 0.03726705645706307
$ws.Range("B6").Value = "[1.0, -0.004467955966576106, -0.013515159629546096, -0.1125505947657036, -0.11787682487752815, 0.08896322718127544, 0.07715368960525949, 0.03842292999415225, 0.03268907777085203, -0.023153673058939457, -0.06142779477611715, -0.028066426470870598, -0.03650418414737522, 0.040067337460878656, 0.0632004838914419, 0.013962545764416896, -0.005557092108048057, -0.05160719264886217, -0.06799564921507192, -0.02244276823969729]"

$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.001351303423638564
$ws.Range("B3").Value = 0.1079204278174231
$ws.Range("B4").Value = 0.00809881279385187
$ws.Range("B5").Value = 0.0000000000006183727995461246
$ws.Range("B7").Value = 0.09155320368891766
$ws.Range("B6").Value = "[0.9999999999999998, 0.002043434080231608, -0.015009377571391298, -0.11122217675688705, -0.1186019290787981, 0.08883855296594748, 0.07866507778417647, 0.037549421691318904, 0.03511525868919956, -0.024665241930410697, -0.06184083780473007, -0.026134983194632214, -0.03662565346173157, 0.04001098834157515, 0.063434271360994, 0.013826161458745195, -0.0048687207625587505, -0.05308640326266957, -0.06779493303737066, -0.020434405969308674]"

$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.002889033725805676
$ws.Range("B3").Value = 0.09577388153613084
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0.00000003685951251672595
$ws.Range("B7").Value = 0.09070738144872696
$ws.Range("B9").Value = 0.7581471104568333
$ws.Range("B6").Value = "[1.0, 0.008607096596604247, -0.05221225668923381, -0.13561830953664733, -0.1168871903398202, 0.08357331039956306, 0.07980028961906957, 0.05977046282477866, 0.020404713490187933, -0.03885714891526289, -0.0544633338689825, -0.02478643711153862, -0.028127775027931498, 0.03724041166277746, 0.0636151010891946, 0.02657692101524521, -0.02314021020126849, -0.054875641863810165, -0.0758410459849212, -0.03570149011423548]"
